$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "VALOR MORA" total value
$ws.Range("E11").Value = 340586

# 2) Update "Cant. Trabajadores" count
$ws.Range("C13").Value = 1

# 3) Update value for first period row (2009) of HAYDIS MIRANDA MARTINEZ
$ws.Range("F16").Value = 35112

# 4) Copy the formatting of the former last row (35, "closing" border style) onto the
#    row that will become the new last row of the (now single-worker) table (row 25),
#    then give it the values that belonged to the old "2106" row (26).
$ws.Range("B35:J35").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)

$ws.Range("E25").Value = "2106"
$ws.Range("F25").Value = 24578
$ws.Range("G25").Value = 877803
$ws.Range("H25").Value = ""
$ws.Range("I25").Value = ""
$ws.Range("J25").Value = ""

# 5) Remove the whole second worker's table (rows 26-35), which also removes the
#    now-superseded original row 35. Rows below shift up automatically.
$ws.Range("B26:J35").EntireRow.Delete()

Write-Output "done"
